# The deck currently uses the "Integral" theme (ppt/theme/theme2.xml) for
# its slide master / slides, while ppt/theme/theme1.xml (unused by any
# slide, only wired to the notes master) carries the stock "Office Theme"
# palette. The authored change swaps the two themes' content so that the
# presentation's live design becomes the default Office Theme palette.
#
# Re-point every theme colour slot on the live (slide-facing) theme to the
# Office Theme RGB values via the Design/ColorScheme object model - this is
# the supported COM surface for rewriting a theme's <a:clrScheme> colours.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order matches MsoThemeColorSchemeIndex: 1=dk1 2=lt1 3=dk2 4=lt2
# 5=accent1 6=accent2 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
$tcs.Item(1).RGB  = 0x00000000   # dk1      000000
$tcs.Item(2).RGB  = 0x00FFFFFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x006A5444   # dk2      44546A
$tcs.Item(4).RGB  = 0x00E6E6E7   # lt2      E7E6E6
$tcs.Item(5).RGB  = 0x00D59B5B   # accent1  5B9BD5
$tcs.Item(6).RGB  = 0x00317DED   # accent2  ED7D31
$tcs.Item(7).RGB  = 0x00A5A5A5   # accent3  A5A5A5
$tcs.Item(8).RGB  = 0x0000C0FF   # accent4  FFC000
$tcs.Item(9).RGB  = 0x00C47244   # accent5  4472C4
$tcs.Item(10).RGB = 0x0047AD70   # accent6  70AD47
$tcs.Item(11).RGB = 0x00C16305   # hlink    0563C1
$tcs.Item(12).RGB = 0x00724F95   # folHlink 954F72
